$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1361
$ws.Range("C2").Value = 1453
$ws.Range("D2").Value = 518
$ws.Range("E2").Value = 376

$ws.Range("B3").Value = 1482
$ws.Range("C3").Value = 1541
$ws.Range("D3").Value = 450
$ws.Range("E3").Value = 349

$ws.Range("H16").Select()
